# Predicted_LaLiga2025_26_table_matchday_0.xlsx
# Adds WIN / TOP4 / TOP5 / TOP6 / RELEGATION placeholder columns (C:G) ahead of
# the existing ExpPoints column (now moved to H), in preparation for a Monte
# Carlo simulation. Also refreshes the Team order / ExpPoints values for matchday 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: insert the new stat columns between Team and ExpPoints ---
$ws.Range("H1").Value = "ExpPoints"   # ExpPoints header moves from C1 to H1
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP4"
$ws.Range("E1").Value = "TOP5"
$ws.Range("F1").Value = "TOP6"
$ws.Range("G1").Value = "RELEGATION"

# Give the new header cells the same style as the existing header row (A1:B1)
$ws.Range("A1").Copy() | Out-Null
$ws.Range("C1:H1").PasteSpecial(-4122) | Out-Null

# --- Data rows: Rank (A) is unchanged; Team (B) is re-ordered; ExpPoints (old C)
# moves to column H; the new C:G stat columns are left blank placeholders for the
# upcoming Monte Carlo simulation output. ---

$ws.Range("B2").Value = "Barcelona"
$ws.Range("H2").Value = 82.43156414117615
$ws.Range("C2:G2").ClearContents()

$ws.Range("B3").Value = "Real Madrid"
$ws.Range("H3").Value = 82.23069947572209
$ws.Range("C3:G3").ClearContents()

$ws.Range("B4").Value = "Atlético de Madrid"
$ws.Range("H4").Value = 74.04488542161154
$ws.Range("C4:G4").ClearContents()

$ws.Range("B5").Value = "Villarreal"
$ws.Range("H5").Value = 63.12052157690897
$ws.Range("C5:G5").ClearContents()

$ws.Range("B6").Value = "Real Betis"
$ws.Range("H6").Value = 58.33862352181674
$ws.Range("C6:G6").ClearContents()

$ws.Range("B7").Value = "Athletic Club"
$ws.Range("H7").Value = 53.29645539727854
$ws.Range("C7:G7").ClearContents()

$ws.Range("B8").Value = "Celta de Vigo"
$ws.Range("H8").Value = 51.11705844163339
$ws.Range("C8:G8").ClearContents()

$ws.Range("B9").Value = "Rayo Vallecano"
$ws.Range("H9").Value = 49.34255759916618
$ws.Range("C9:G9").ClearContents()

$ws.Range("B10").Value = "Osasuna"
$ws.Range("H10").Value = 46.88083555078431
$ws.Range("C10:G10").ClearContents()

$ws.Range("B11").Value = "Real Sociedad"
$ws.Range("H11").Value = 46.27098356224886
$ws.Range("C11:G11").ClearContents()

$ws.Range("B12").Value = "Espanyol"
$ws.Range("H12").Value = 44.81054064960833
$ws.Range("C12:G12").ClearContents()

$ws.Range("B13").Value = "Valencia"
$ws.Range("H13").Value = 44.79671429353233
$ws.Range("C13:G13").ClearContents()

$ws.Range("B14").Value = "Sevilla"
$ws.Range("H14").Value = 44.75928889404506
$ws.Range("C14:G14").ClearContents()

$ws.Range("B15").Value = "Alavés"
$ws.Range("H15").Value = 43.91333142750105
$ws.Range("C15:G15").ClearContents()

$ws.Range("B16").Value = "Getafe"
$ws.Range("H16").Value = 43.61816622394461
$ws.Range("C16:G16").ClearContents()

$ws.Range("B17").Value = "Mallorca"
$ws.Range("H17").Value = 43.26290762407913
$ws.Range("C17:G17").ClearContents()

$ws.Range("B18").Value = "Elche"
$ws.Range("H18").Value = 41.66700788472858
$ws.Range("C18:G18").ClearContents()

$ws.Range("B19").Value = "Levante"
$ws.Range("H19").Value = 39.2591493820248
$ws.Range("C19:G19").ClearContents()

$ws.Range("B20").Value = "Girona"
$ws.Range("H20").Value = 36.42404208697945
$ws.Range("C20:G20").ClearContents()

$ws.Range("B21").Value = "Real Oviedo"
$ws.Range("H21").Value = 34.68837316042448
$ws.Range("C21:G21").ClearContents()
